# Update "想去人数" (want-to-go count) values in column F
# Sheet "展览" (sheet index 1) and sheet "全部类型" (sheet index 4)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value  = 862
$ws1.Range("F8").Value  = 4756
$ws1.Range("F9").Value  = 4756
$ws1.Range("F12").Value = 164
$ws1.Range("F16").Value = 7601
$ws1.Range("F21").Value = 535
$ws1.Range("F22").Value = 1402
$ws1.Range("F24").Value = 6288
$ws1.Range("F28").Value = 6198
$ws1.Range("F33").Value = 448
$ws1.Range("F46").Value = 453

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F9").Value  = 4756
$ws4.Range("F10").Value = 4756
$ws4.Range("F13").Value = 164
$ws4.Range("F17").Value = 7601
$ws4.Range("F20").Value = 535
$ws4.Range("F21").Value = 1402
$ws4.Range("F23").Value = 6288
$ws4.Range("F29").Value = 6198
$ws4.Range("F46").Value = 453
